$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.166.32"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.791.07"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'226.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'31.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").Value = "'0.0690"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("D11").Value = "'0.0946"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "2.048.26"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").Value = "'11.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "1.787.65"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "34.101.05"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").Value = "'68.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").Value = "'245.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'10.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").Value = "'4.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "'2.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").Value = "'161.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "'7.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").Value = "'16.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "'0.114"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "'3.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.15%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "1.456.64"
$ws.Range("E35").Value = "  +4.53%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").Value = "'2.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.78%  "
$ws.Range("E38").Value = "  +2.59%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "'80.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.35%  "
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").Value = "'0.922"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("D43").Value = "'2.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "0.0₆0136"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").Value = "1.949.85"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").Value = "'106.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("E51").Value = "  +0.00%  "
